$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: cohort 2024, period 1 - num_customers 119 -> 120, retention_rate recalculated
$ws.Range("C36").Value = 120
$ws.Range("E36").Value = 120 / 1930

# Row 37: cohort 2025, period 0 - num_customers 739 -> 750, cohort_size 739 -> 750
$ws.Range("C37").Value = 750
$ws.Range("D37").Value = 750
